$d = $word.ActiveDocument

# Locate the anchor paragraph: the last question in the bulleted list,
# ending with "3242XZY78K3". New "Regle de gestion" content is inserted
# right after it and before the following table.
$anchorIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
  if ($d.Paragraphs($i).Range.Text -like "*3242XZY78K3*") {
    $anchorIndex = $i
    break
  }
}
if ($anchorIndex -eq -1) {
  throw "Anchor paragraph not found"
}

$fragments = @(
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="708"/><w:rPr><w:color w:val="FF0000"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr><w:r><w:rPr><w:color w:val="FF0000"/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:lastRenderedPageBreak/><w:t>Règle de gestion :</w:t></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="708"/></w:pPr><w:r><w:t>Un avion appartient à un propriétaire</w:t></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="708"/></w:pPr><w:r><w:t>Un propriétaire à zéro ou plusieurs avions</w:t></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="708"/></w:pPr></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="708"/></w:pPr><w:r><w:t>Un avion à un type</w:t></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="708"/></w:pPr><w:r><w:t>Un type correspond à un ou plusieurs avions</w:t></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="708"/></w:pPr></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="708"/></w:pPr><w:r><w:t>Un avion peut avoir une intervention</w:t></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="708"/></w:pPr><w:r><w:t>Une intervention correspond à un avion</w:t></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="708"/></w:pPr></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="708"/></w:pPr><w:r><w:t>Un avion est piloté par un pilote</w:t></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="708"/></w:pPr><w:r><w:t xml:space="preserve">Un pilote peut </w:t></w:r><w:r><w:t>piloter</w:t></w:r><w:r><w:t xml:space="preserve"> plusieurs avions</w:t></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="708"/></w:pPr></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="708"/></w:pPr><w:r><w:t>Un propriétaire peut avoir zéro ou plusieurs types d’avion</w:t></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="708"/></w:pPr><w:r><w:t>Un type d’avion peut être détenu par zéro ou plusieurs propriétaires</w:t></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="708"/></w:pPr></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="708"/></w:pPr><w:r><w:t>Un type d’avion est entretenu par deux ou plusieurs mécaniciens</w:t></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="708"/></w:pPr><w:r><w:t xml:space="preserve">Un mécanicien entretient un ou </w:t></w:r><w:r><w:t>plusieurs types</w:t></w:r><w:r><w:t xml:space="preserve"> d’avions</w:t></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="708"/></w:pPr></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="708"/></w:pPr><w:r><w:t>Un type d’avion peut être piloté par zéro ou plusieurs pilotes</w:t></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="708"/></w:pPr><w:r><w:t>Un pilote peut piloter un ou plusieurs types d’avions</w:t></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="708"/></w:pPr></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="708"/></w:pPr><w:r><w:t>Un mécanicien peut effectuer à zéro ou plusieurs interventions</w:t></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="708"/></w:pPr><w:r><w:t>Une intervention est effectuée par au moins deux mécaniciens</w:t></w:r></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="708"/></w:pPr></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="708"/></w:pPr></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="708"/></w:pPr></w:p>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>',
  '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:ind w:firstLine="708"/></w:pPr></w:p>'
)

$cursor = $anchorIndex
foreach ($frag in $fragments) {
  $anchorPara = $d.Paragraphs($cursor)
  $null = $anchorPara.Range.InsertParagraphAfter()
  $cursor = $cursor + 1
  $newPara = $d.Paragraphs($cursor)
  $null = $newPara.Range.InsertXML($frag)
}

Write-Output "Inserted $($fragments.Count) paragraphs after index $anchorIndex. New paragraph count: $($d.Paragraphs.Count)"
